$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the formatting of the last existing data row (row 80) down into the
# new row 81, then overwrite it with the data for Post71.
$ws.Range("B80:F80").Copy($ws.Range("B81:F81"))

$ws.Cells.Item(81, 2).Value = 71
$ws.Cells.Item(81, 3).Value = "Variable Partitioning | Operating System - M05 P04"
$ws.Cells.Item(81, 4).Value = 44183
$ws.Cells.Item(81, 5).Value = "https://programmingport.hashnode.dev/variable-partitioning-or-operating-system-m05-p04"
$ws.Cells.Item(81, 6).Value = "https://dev.to/rahulmishra05/variable-partitioning-operating-system-m05-p04-3g9a"

$excel.CutCopyMode = $false

# Expand the table (ListObject) to include the new row
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B10:F81"))

# Update the active selection to match the new last cell
$ws.Range("E81").Select() | Out-Null
